$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 used to be "Transistors (x4)/BC337" with a "to test" comment.
# The transistor line is removed; the diode line (formerly row 10) moves
# up into row 9 and its comment is resolved to "Ok".
$ws.Range("A9").Value = "Diodes (x4)"
$ws.Range("B9").Value = "1N4001"
$ws.Range("C9").Formula = "=0.21/10*4"
$ws.Range("E9").Value = "Ok"

# Row 10 becomes the new MOSFETs (x4) line that was just ordered.
$ws.Range("A10").Value = "MOSFETs (x4)"
$ws.Range("B10").Value = "PMV31XN"
$ws.Range("C10").Formula = "=0.286*4"
$ws.Range("E10").Value = "Ordered"
$ws.Range("D10").Value = "RS Online"

# The pager motors and props have now been received.
$ws.Range("E11").Value = "Ok"
$ws.Range("E12").Value = "Ok"

# Leave the active selection on the newly edited cell, as in the saved file.
$ws.Range("D10").Select()
